$p = $ppt.ActivePresentation

# --- Update the three result tables to use the new table style ---
$newStyleId = "{38AC001A-CA69-435E-AB2D-45CE6F1DD35A}"

foreach ($slideIdx in 9, 10, 13) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- Simplify the "We beat State of the Art..." bullet on slide 11 ---
$s11 = $p.Slides.Item(11)
$bodyShape = $s11.Shapes.Item(2)
$tr = $bodyShape.TextFrame.TextRange
$para2 = $tr.Paragraphs(2)
$chars = $tr.Characters($para2.Start, $para2.Length)
$chars.Text = "We beat State of the Art accuracy (0.274) using pre-trained RoBERTa (0.292) model."
